$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# Fill in priority values (column B) for the leaf task rows
$ws.Range("B4").Value = 3    # Regisztráció
$ws.Range("B5").Value = 3    # Belépés
$ws.Range("B7").Value = 1    # Feltöltés
$ws.Range("B8").Value = 1    # Letöltés
$ws.Range("B10").Value = 2   # Listázás
$ws.Range("B11").Value = 2   # Törlés
$ws.Range("B12").Value = 2   # Módosítás
$ws.Range("B13").Value = 2   # Keresés
$ws.Range("B14").Value = 2   # CSS formázás
$ws.Range("B15").Value = 1   # Tesztelés
$ws.Range("B16").Value = 3   # Logolás

# Update the saved selection to B1
$ws.Range("B1").Select()
